$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2-7 from 45170 to 45174
foreach ($r in 2..7) {
    $ws.Cells.Item($r, 3).Value = 45174
}
